$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: account holder name / card number
$ws.Range("C2").Value = "Hartmut"

# B3 holds a 16-digit card number that must stay text (Excel would
# otherwise auto-coerce a pure-digit string to a number). Force text via
# NumberFormat, then copy the original cell's formatting back on top so
# the stored style index is unchanged.
$b3 = $ws.Range("B3")
$b3.NumberFormat = "@"
$b3.Value = "2570314725427075"
$ws.Range("C3").Copy()
$b3.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("C3").Value = "Mohaupt"

# Row 5: opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 17.06.2024"

# Row 6
$ws.Range("B6").Value = "20.06."
$ws.Range("C6").Value = "21.06."
$ws.Range("D6").Value = "AMAZON.DE MKTPLC EU ORONXD"
$ws.Range("E6").Value = "191,56-"

# Row 7
$ws.Range("B7").Value = "22.06."
$ws.Range("C7").Value = "23.06."
$ws.Range("D7").Value = "AMAZON.DE MKTPLC EU GDNQXS"
$ws.Range("E7").Value = "78,33-"

# Row 8
$ws.Range("B8").Value = "25.06."
$ws.Range("C8").Value = "26.06."
$ws.Range("D8").Value = "RECHNUNG VODAFONE GMBH 78336771"
$ws.Range("E8").Value = "41,30-"

# Row 9
$ws.Range("B9").Value = "26.06."
$ws.Range("C9").Value = "27.06."
$ws.Range("D9").Value = "PAYPAL UKMOQE"
$ws.Range("E9").Value = "48,56-"

# Row 10
$ws.Range("B10").Value = "27.06."
$ws.Range("C10").Value = "28.06."
$ws.Range("D10").Value = "BEITRAG Allianz SE K-82053496"
$ws.Range("E10").Value = "55,83-"

# Row 11: the 5th transaction is removed entirely - clear all four cells.
# Note: Range.ClearContents() misbehaves on merged cells (E11 is part of
# the E11:F11 merge) in this runtime, so assign empty strings directly.
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$e11 = $ws.Range("E11")
$e11.Value = ""
$e11.HorizontalAlignment = -4152  # xlRight
$e11.VerticalAlignment = -4108    # xlCenter
$e11.WrapText = $true

# Row 12: closing balance
$ws.Range("D12").Value = "KONTOSTAND AM 30.06.2024"
$ws.Range("E12").Value = "415,58-"

# Row 13: next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 06.07.2024"
